$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "legenda"
$ws.Range("B1").Value = "area"
$ws.Range("D1").Value = "area_km2"

$ws.Range("B2").Value = 72547.90897610001
$ws.Range("B3").Value = 86240.7412355
